# KIBON-122: Angebotstyp auf GUI und in Reports hinzugefuegt
#
# Adds a new "Betreuungsangebot Typ" column to the ZahlungAuftragPeriode
# report sheet, between the existing "Institution" and "Gemeinde" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new column at B. This shifts every existing cell in
# columns B..E one column to the right (B->C, C->D, D->E, E->F) and also
# carries the column-width metadata along with it (old col B's custom
# width moves to col C, old col C/D/E's auto widths move to col D/E/F,
# and the default-width block that used to start at col 6 now starts at
# col 7) - exactly matching how the header/data rows (6 and 7) and the
# decorative empty rows (1, 2, 3) gained a new column.
$ws.Columns.Item(2).Insert()

# Row 4 ({periodeTitle}/{periode}) is not part of the table and must keep
# its original 2-cell layout (A4, B4) - the blanket column insert above
# shifted its B4 value into C4, so move it back and drop the now-empty
# C4 cell entirely so the row goes back to only spanning A:B.
$ws.Range("C4").Copy($ws.Range("B4"))
$ws.Range("C4").Clear()

# The newly inserted column B is currently an empty clone of the old
# column B (institution/gemeinde) in rows 6 and 7; give it its own
# header and placeholder text.
$ws.Range("B6").Value = "{betreuungsangebotTypTitle}"
$ws.Range("B7").Value = "{betreuungsangebotTyp}"

# Match the author's explicit width for the freshly inserted column.
$ws.Range("B1").EntireColumn.ColumnWidth = 15

Write-Host "Added Betreuungsangebot Typ column to report"
